$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that was bumped by one day
# (45181 -> 45182) for every data row (rows 2 through 532).
$ws.Range("C2:C532").Value = 45182
